$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.618.07'
$ws.Range("E2").Value = '  -2.59%  '

$ws.Range("D3").Value = '1.666.37'
$ws.Range("E3").Value = '  -3.76%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.56'
$ws.Range("E5").Value = '  -2.31%  '

$ws.Range("E6").Value = '  -2.19%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.06'
$ws.Range("E8").Value = '  -1.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.264'
$ws.Range("E9").Value = '  -1.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0621'
$ws.Range("E10").Value = '  -2.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0880'
$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("D12").Value = '1.901.43'
$ws.Range("E12").Value = '  -3.74%  '

$ws.Range("D13").Value = '1.667.93'
$ws.Range("E13").Value = '  -3.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.14'
$ws.Range("E14").Value = '  -3.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.566'
$ws.Range("E15").Value = '  +0.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.43'
$ws.Range("E16").Value = '  -1.98%  '

$ws.Range("D17").Value = '27.603.04'
$ws.Range("E17").Value = '  -2.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.49'
$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").Value = '0.0₃0731'
$ws.Range("E19").Value = '  -3.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.69'
$ws.Range("E20").Value = '  -4.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.50'
$ws.Range("E22").Value = '  -3.28%  '

$ws.Range("E23").Value = '  -3.62%  '

$ws.Range("E24").Value = '  -3.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.59'
$ws.Range("E25").Value = '  -1.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.24'
$ws.Range("E26").Value = '  -4.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.45'
$ws.Range("E27").Value = '  -1.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("E29").Value = '  -2.39%  '

$ws.Range("E30").Value = '  +1.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0504'
$ws.Range("E31").Value = '  -1.82%  '

$ws.Range("E32").Value = '  -2.65%  '

$ws.Range("D33").Value = '1.460.47'
$ws.Range("E33").Value = '  -3.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.13'
$ws.Range("E34").Value = '  -4.55%  '

$ws.Range("E35").Value = '  -4.96%  '

$ws.Range("E36").Value = '  -1.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.928'
$ws.Range("E37").Value = '  -4.31%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.578'
$ws.Range("E38").Value = '  -5.05%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0173'
$ws.Range("E39").Value = '  -1.97%  '

$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '69.67'
$ws.Range("E40").Value = '  -1.79%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.03'
$ws.Range("E41").Value = '  -3.92%  '

$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("E43").Value = '  -3.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.40'
$ws.Range("E44").Value = '  -5.84%  '

$ws.Range("D45").Value = '1.809.56'
$ws.Range("E45").Value = '  -3.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.791'
$ws.Range("E46").Value = '  -1.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.74'
$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("E48").Value = '  -2.39%  '

$ws.Range("E49").Value = '  -5.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.103'
$ws.Range("E50").Value = '  -2.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.93'
$ws.Range("E51").Value = '  -3.41%  '
